$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.908.80'
$ws.Range("E2").Value = '  -0.35%  '

$ws.Range("D3").Value = '1.622.21'
$ws.Range("E3").Value = '  -1.18%  '

$ws.Range("E4").Value = '  +0.26%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.50'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.07%  '

$ws.Range("E6").Value = '  -0.97%  '

$ws.Range("E7").Value = '  +0.25%  '

$ws.Range("E8").Value = '  -2.55%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0614'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.81%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.22'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -6.62%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0786'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.29%  '

$ws.Range("D12").Value = '1.847.78'
$ws.Range("E12").Value = '  -1.14%  '

$ws.Range("D13").Value = '1.618.98'
$ws.Range("E13").Value = '  -1.23%  '

$ws.Range("E14").Value = '  -2.33%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.522'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -3.87%  '

$ws.Range("D16").Value = '25.899.69'
$ws.Range("E16").Value = '  -0.52%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.06'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.64%  '

$ws.Range("D18").Value = '0.0₃0731'
$ws.Range("E18").Value = '  -4.18%  '

$ws.Range("E19").Value = '  +0.24%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '191.70'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.18%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.22'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.25%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.54'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -3.82%  '

$ws.Range("E23").Value = '  -2.34%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.132'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.62%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.70'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.53%  '

$ws.Range("E26").Value = '  +0.36%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.74'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -3.21%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.70'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.48%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.12'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.58%  '

$ws.Range("E30").Value = '  -1.50%  '

$ws.Range("E31").Value = '  -2.77%  '

$ws.Range("E32").Value = '  -4.65%  '

$ws.Range("E33").Value = '  -5.70%  '

$ws.Range("E34").Value = '  -3.05%  '

$ws.Range("E35").Value = '  -2.41%  '

$ws.Range("D36").Value = '1.114.66'

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.844'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -6.66%  '

$ws.Range("E38").Value = '  -1.27%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.515'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -4.55%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0152'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.69%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '97.88'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.13%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.764'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.11%  '

$ws.Range("D43").Value = '1.758.59'
$ws.Range("E43").Value = '  -1.10%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.15'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -5.84%  '

$ws.Range("E45").Value = '  -1.61%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0529'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.50%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.20'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -4.16%  '

$ws.Range("E48").Value = '  -1.80%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.413'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.35%  '

$ws.Range("E50").Value = '  +0.44%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.46'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -3.49%  '
